$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Add the new data for row 42 (Responsable + completed flag + link to new article)
$ws.Range("B42").Value = "Agustina"
$ws.Range("C42").Value = 1
$ws.Range("C42").NumberFormat = "0%"
$ws.Range("D42").Value = "http://mygnet.net/articulos/java/creacion_de_graficos_en_ireport.707"

# Update the active selection to reflect where the user ended up editing
$ws.Activate()
$ws.Range("D43").Select()
